$d = $word.ActiveDocument

$replacements = @(
    @{old="702÷7="; new="159÷5="},
    @{old="808÷4="; new="464÷8="},
    @{old="333÷9="; new="687÷3="},
    @{old="821÷9="; new="573÷7="},
    @{old="379÷8="; new="334÷3="},
    @{old="959÷2="; new="646÷3="},
    @{old="396÷2="; new="878÷8="},
    @{old="314÷2="; new="823÷8="},
    @{old="731÷3="; new="666÷4="},
    @{old="443÷7="; new="366÷8="},
    @{old="629÷7="; new="828÷6="},
    @{old="900÷6="; new="209÷4="},
    @{old="702÷8="; new="971÷2="},
    @{old="509÷5="; new="161÷9="},
    @{old="582÷6="; new="982÷7="},
    @{old="699÷2="; new="589÷7="},
    @{old="242÷6="; new="158÷6="},
    @{old="810÷6="; new="705÷8="},
    @{old="872÷9="; new="490÷7="},
    @{old="165÷5="; new="707÷7="},
    @{old="283÷2="; new="141÷5="},
    @{old="668÷5="; new="923÷5="},
    @{old="188÷2="; new="686÷3="},
    @{old="684÷2="; new="676÷5="},
    @{old="472÷9="; new="627÷9="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
